$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ADC count values for rows 2-42 (column B), corresponding to temperatures
# -45..155 in 5-degree steps. (commit: "with excel df - 2048")
$newValues = @(
    285441,
    291690,
    297917,
    304137,
    310564,
    316736,
    322881,
    329856,
    335978,
    342071,
    348949,
    355302,
    361879,
    368191,
    374728,
    381001,
    387790,
    394240,
    400466,
    407200,
    413893,
    420569,
    426415,
    433060,
    439681,
    446274,
    452845,
    459855,
    466392,
    472540,
    479032,
    485934,
    492378,
    499206,
    506013,
    512390,
    518751,
    525455,
    532539,
    539199,
    545409
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("B$row").Value = $newValues[$i]
}
